$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as TEXT, matching the
# original inlineStr cell type, without permanently altering the cells style.
function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range("D2").Value = '70.267.91'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '3.614.21'
$ws.Range("E3").Value = '  +2.21%  '
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.10%  '
Set-TextValue $ws.Range("D5") '602.12'
$ws.Range("E5").Value = '  -0.33%  '
Set-TextValue $ws.Range("D6") '195.66'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("E8").Value = '  +0.06%  '
Set-TextValue $ws.Range("D9") '0.211'
$ws.Range("E9").Value = '  +2.96%  '
Set-TextValue $ws.Range("D10") '0.647'
$ws.Range("E10").Value = '  -0.66%  '
Set-TextValue $ws.Range("D11") '53.22'
$ws.Range("E11").Value = '  -0.99%  '
Set-TextValue $ws.Range("D12") '0.0000304'
$ws.Range("E12").Value = '  +0.40%  '
Set-TextValue $ws.Range("D13") '9.58'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").Value = '4.192.60'
$ws.Range("E14").Value = '  +2.23%  '
Set-TextValue $ws.Range("D15") '603.40'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("E16").Value = '  +1.90%  '
$ws.Range("D17").Value = '70.417.22'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").Value = '3.607.55'
$ws.Range("E18").Value = '  +1.74%  '
Set-TextValue $ws.Range("D19") '19.05'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("E20").Value = '  +1.59%  '
Set-TextValue $ws.Range("D21") '1.00'
$ws.Range("E21").Value = '  +0.72%  '
Set-TextValue $ws.Range("D22") '18.50'
$ws.Range("E22").Value = '  +1.80%  '
Set-TextValue $ws.Range("D23") '5.20'
$ws.Range("E23").Value = '  -0.31%  '
Set-TextValue $ws.Range("D24") '103.37'
$ws.Range("E24").Value = '  +0.75%  '
Set-TextValue $ws.Range("D25") '4.62'
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -4.72%  '
Set-TextValue $ws.Range("D27") '10.59'
$ws.Range("E27").Value = '  -2.90%  '
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("E29").Value = '  +0.85%  '
Set-TextValue $ws.Range("D30") '4.67'
$ws.Range("E30").Value = '  +7.52%  '
$ws.Range("E31").Value = '  +2.04%  '
Set-TextValue $ws.Range("D32") '12.30'
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("E33").Value = '  +2.06%  '
Set-TextValue $ws.Range("D34") '63.24'
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").Value = '0.0₃0882'
$ws.Range("E35").Value = '  +2.64%  '
$ws.Range("D36").Value = '3.934.53'
$ws.Range("E36").Value = '  +5.38%  '
Set-TextValue $ws.Range("D37") '530.72'
$ws.Range("E37").Value = '  +8.62%  '
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("E39").Value = '  +0.10%  '
Set-TextValue $ws.Range("D40") '36.87'
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("E41").Value = '  -0.99%  '
Set-TextValue $ws.Range("D42") '3.53'
$ws.Range("E42").Value = '  -2.94%  '
$ws.Range("E43").Value = '  +0.60%  '
Set-TextValue $ws.Range("D44") '0.0460'
$ws.Range("E44").Value = '  +0.13%  '
Set-TextValue $ws.Range("D45") '3.60'
$ws.Range("E45").Value = '  +9.22%  '
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("E48").Value = '  -0.17%  '
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("E50").Value = '  -1.74%  '
Set-TextValue $ws.Range("D51") '1.29'
$ws.Range("E51").Value = '  +1.03%  '
